# Remove the deprecated "NanoDESI" entry from the assay_type list, leaving
# only "NanoPOTS". This shifts "NanoPOTS" from A2 up to A1 on the
# "assay_type list" sheet, and Excel cleans up the now-unused "NanoDESI"
# shared string automatically on save (re-indexing every other shared
# string reference throughout the workbook).

$wb = $excel.ActiveWorkbook

$assayTypeList = $wb.Worksheets.Item("assay_type list")
$assayTypeList.Range("A1").EntireRow.Delete()

# Update the data validation on the main sheet's assay_type column (L) so
# its list range and error message reflect the single remaining entry.
$mainSheet = $wb.Worksheets.Item("Export as TSV")
$dv = $mainSheet.Range("L2:L1048576").Validation
$dv.Modify(3, 1, 1, "='assay_type list'!`$A`$1:`$A`$1")
$dv.ErrorMessage = "Value must be one of: NanoPOTS."
